$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 2.2

# Row 4
$ws.Range("G4").Value = 3.8
$ws.Range("H4").Value = 3.65
$ws.Range("I4").Value = 1.82
$ws.Range("J4").Value = 4.15
$ws.Range("K4").Value = 2.2
$ws.Range("L4").Value = 2.35
$ws.Range("O4").Value = 1.22
$ws.Range("P4").Value = 3.4
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 1.98
$ws.Range("U4").Value = 1.62
$ws.Range("V4").Value = 2.02
$ws.Range("W4").Value = 12.5
$ws.Range("X4").Value = 22
$ws.Range("Y4").Value = 13
$ws.Range("Z4").Value = 55
$ws.Range("AA4").Value = 32
$ws.Range("AB4").Value = 35
$ws.Range("AC4").Value = 12
$ws.Range("AD4").Value = 7.2
$ws.Range("AE4").Value = 13.5
$ws.Range("AF4").Value = 55
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 8.25
$ws.Range("AI4").Value = 9.5
$ws.Range("AJ4").Value = 8.25
$ws.Range("AK4").Value = 15.5
$ws.Range("AL4").Value = 13.5
$ws.Range("AM4").Value = 22
$ws.Range("AN4").Value = 5.7
$ws.Range("AO4").Value = 21
$ws.Range("AP4").Value = 26
$ws.Range("AQ4").Value = 110
$ws.Range("AR4").Value = 150
$ws.Range("AS4").Value = 300
$ws.Range("AT4").Value = 2.9
$ws.Range("AU4").Value = 7
$ws.Range("AV4").Value = 55
$ws.Range("AX4").Value = 3.75
$ws.Range("AY4").Value = 8.75
$ws.Range("AZ4").Value = 16.5
$ws.Range("BA4").Value = 30
$ws.Range("BB4").Value = 55
$ws.Range("BC4").Value = 200

# Row 5
$ws.Range("O5").Value = 1.33
$ws.Range("P5").Value = 3.4
$ws.Range("Q5").Value = 2.07
$ws.Range("R5").Value = 1.83

# Row 9
$ws.Range("G9").Value = 1.53
$ws.Range("H9").Value = 3.8
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 2.2
$ws.Range("K9").Value = 2.1
$ws.Range("L9").Value = 7
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("X9").Value = 6
$ws.Range("Z9").Value = 10
$ws.Range("AD9").Value = 7.5
$ws.Range("AE9").Value = 23
$ws.Range("AH9").Value = 13
$ws.Range("AI9").Value = 34
$ws.Range("AK9").Value = 81
$ws.Range("AM9").Value = 67
$ws.Range("AN9").Value = 3.25
$ws.Range("AO9").Value = 8
$ws.Range("AQ9").Value = 26
$ws.Range("AX9").Value = 8
$ws.Range("AY9").Value = 41

# Row 12
$ws.Range("M12").Value = 1.04
$ws.Range("N12").Value = 13
$ws.Range("Q12").Value = 1.83
$ws.Range("R12").Value = 2.03

# Row 13
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 15

# Row 16
$ws.Range("G16").Value = 1.75
$ws.Range("H16").Value = 3.3
$ws.Range("I16").Value = 4.5
$ws.Range("J16").Value = 2.4
$ws.Range("K16").Value = 2.1
$ws.Range("L16").Value = 5
$ws.Range("M16").Value = 1.07
$ws.Range("N16").Value = 9
$ws.Range("O16").Value = 1.33
$ws.Range("P16").Value = 3.25
$ws.Range("Q16").Value = 2.05
$ws.Range("R16").Value = 1.75
$ws.Range("S16").Value = 1.44
$ws.Range("T16").Value = 2.63
$ws.Range("U16").Value = 1.91
$ws.Range("V16").Value = 1.8
$ws.Range("W16").Value = 6.5
$ws.Range("X16").Value = 7.5
$ws.Range("Y16").Value = 8.5
$ws.Range("Z16").Value = 13
$ws.Range("AA16").Value = 15
$ws.Range("AB16").Value = 29
$ws.Range("AC16").Value = 9
$ws.Range("AD16").Value = 6.5
$ws.Range("AE16").Value = 17
$ws.Range("AF16").Value = 51
$ws.Range("AG16").Value = 351
$ws.Range("AH16").Value = 12
$ws.Range("AI16").Value = 23
$ws.Range("AJ16").Value = 17
$ws.Range("AK16").Value = 51
$ws.Range("AL16").Value = 41
$ws.Range("AM16").Value = 41
$ws.Range("AN16").Value = 3.6
$ws.Range("AO16").Value = 9.5
$ws.Range("AP16").Value = 21
$ws.Range("AQ16").Value = 34
$ws.Range("AR16").Value = 51
$ws.Range("AS16").Value = 151
$ws.Range("AT16").Value = 2.63
$ws.Range("AU16").Value = 8.5
$ws.Range("AV16").Value = 67
$ws.Range("AX16").Value = 6.5
$ws.Range("AY16").Value = 26
$ws.Range("AZ16").Value = 34
$ws.Range("BA16").Value = 101
$ws.Range("BB16").Value = 126
$ws.Range("BC16").Value = 301

# Row 17
$ws.Range("G17").Value = 1.17
$ws.Range("H17").Value = 6.5
$ws.Range("I17").Value = 11
$ws.Range("K17").Value = 2.88
$ws.Range("L17").Value = 11
$ws.Range("N17").Value = 17
$ws.Range("O17").Value = 1.14
$ws.Range("P17").Value = 5.5
$ws.Range("R17").Value = 2.5
$ws.Range("S17").Value = 1.25
$ws.Range("T17").Value = 3.75
$ws.Range("U17").Value = 2.25
$ws.Range("V17").Value = 1.57
$ws.Range("W17").Value = 8
$ws.Range("X17").Value = 6
$ws.Range("Y17").Value = 10
$ws.Range("Z17").Value = 7
$ws.Range("AA17").Value = 11
$ws.Range("AB17").Value = 34
$ws.Range("AC17").Value = 17
$ws.Range("AD17").Value = 13
$ws.Range("AE17").Value = 29
$ws.Range("AF17").Value = 81
$ws.Range("AH17").Value = 26
$ws.Range("AI17").Value = 51
$ws.Range("AJ17").Value = 34
$ws.Range("AK17").Value = 151
$ws.Range("AL17").Value = 81
$ws.Range("AM17").Value = 67
$ws.Range("AN17").Value = 3.1
$ws.Range("AO17").Value = 5
$ws.Range("AP17").Value = 17
$ws.Range("AQ17").Value = 12
$ws.Range("AR17").Value = 41
$ws.Range("AS17").Value = 126
$ws.Range("AT17").Value = 3.75
$ws.Range("AU17").Value = 11
$ws.Range("AV17").Value = 67
$ws.Range("AY17").Value = 51
$ws.Range("AZ17").Value = 51
$ws.Range("BA17").Value = 351
$ws.Range("BB17").Value = 301

# Row 18
$ws.Range("G18").Value = 3
$ws.Range("I18").Value = 2.45
$ws.Range("J18").Value = 3.6
$ws.Range("S18").Value = 1.44
$ws.Range("T18").Value = 2.63
$ws.Range("AB18").Value = 34
$ws.Range("AC18").Value = 9.5
$ws.Range("AF18").Value = 51
$ws.Range("AH18").Value = 8
$ws.Range("AJ18").Value = 9.5
$ws.Range("AL18").Value = 19
$ws.Range("AO18").Value = 17
$ws.Range("AP18").Value = 26
$ws.Range("AR18").Value = 81
$ws.Range("AS18").Value = 201
$ws.Range("AT18").Value = 2.63

# Row 20
$ws.Range("G20").Value = 2.77
$ws.Range("H20").Value = 3.25
$ws.Range("I20").Value = 2.42
$ws.Range("K20").Value = 2.05
$ws.Range("L20").Value = 3.05
$ws.Range("M20").Value = 1.07
$ws.Range("N20").Value = 7
$ws.Range("O20").Value = 1.32
$ws.Range("P20").Value = 3.1
$ws.Range("Q20").Value = 1.93
$ws.Range("R20").Value = 1.78
$ws.Range("S20").Value = 1.44
$ws.Range("T20").Value = 2.6
$ws.Range("U20").Value = 1.75
$ws.Range("V20").Value = 1.98
$ws.Range("W20").Value = 8.5
$ws.Range("X20").Value = 14
$ws.Range("Z20").Value = 32
$ws.Range("AC20").Value = 7
$ws.Range("AD20").Value = 6.2
$ws.Range("AF20").Value = 65
$ws.Range("AG20").Value = 500
$ws.Range("AH20").Value = 8.25
$ws.Range("AI20").Value = 12
$ws.Range("AK20").Value = 26
$ws.Range("AL20").Value = 19.5
$ws.Range("AM20").Value = 29
$ws.Range("AN20").Value = 4.7
$ws.Range("AP20").Value = 24
$ws.Range("AQ20").Value = 75
$ws.Range("AR20").Value = 120
$ws.Range("AS20").Value = 350
$ws.Range("AT20").Value = 2.6
$ws.Range("AU20").Value = 7.2
$ws.Range("AV20").Value = 70
$ws.Range("AW20").Value = 151
$ws.Range("AX20").Value = 4.3
$ws.Range("AY20").Value = 13

# Row 25
$ws.Range("G25").Value = 1.27
$ws.Range("H25").Value = 5.2
$ws.Range("I25").Value = 9.25
$ws.Range("J25").Value = 1.7
$ws.Range("K25").Value = 2.55
$ws.Range("L25").Value = 7.4
$ws.Range("M25").Value = 1.02
$ws.Range("N25").Value = 17
$ws.Range("O25").Value = 1.18
$ws.Range("P25").Value = 4.85
$ws.Range("Q25").Value = 1.53
$ws.Range("R25").Value = 2.2
$ws.Range("S25").Value = 1.28
$ws.Range("T25").Value = 3.48
$ws.Range("W25").Value = 7.6
$ws.Range("X25").Value = 6.3
$ws.Range("Y25").Value = 9
$ws.Range("Z25").Value = 7.6
$ws.Range("AA25").Value = 10.75
$ws.Range("AC25").Value = 14.5
$ws.Range("AD25").Value = 10.75
$ws.Range("AE25").Value = 24
$ws.Range("AH25").Value = 25
$ws.Range("AI25").Value = 70
$ws.Range("AJ25").Value = 29
$ws.Range("AK25").Value = 250
$ws.Range("AL25").Value = 110
$ws.Range("AM25").Value = 90
$ws.Range("AO25").Value = 5.4
$ws.Range("AP25").Value = 15.5
$ws.Range("AQ25").Value = 13.5
$ws.Range("AR25").Value = 40
$ws.Range("AT25").Value = 3.3
$ws.Range("AU25").Value = 9
$ws.Range("AX25").Value = 9.75
$ws.Range("AY25").Value = 50
$ws.Range("BA25").Value = 400
$ws.Range("BB25").Value = 350
